$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Leche La Lechera"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 200
